$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AM (39) previously shared the default width (11.5546875) with columns
# AH:AN (34-42). It now gets its own, wider column definition.
$ws.Columns.Item(39).ColumnWidth = 15

# Rows 2-39, 41-42, 45-47, 49-50 had an empty AM cell; they now hold the value 0
# (rows 40, 43, 44, 48 and 51 already contained numbers and stay untouched).
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,41,42,45,46,47,49,50)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 39).Value = 0
}

# Update the view/selection state to match the saved workbook.
$ws.Activate()
$ws.Range("AM50").Select()
